$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume table update (GitHub Actions refresh)
$updates = @(
    @{Row=2; D='28.701.35'; E='  +2.52%  '},
    @{Row=3; D='1.906.16'; E='  +2.64%  '},
    @{Row=4; D='1.036'; E='  +3.13%  '},
    @{Row=5; D='320.18'; E='  +2.55%  '},
    @{Row=6; D='1.030'; E='  +2.69%  '},
    @{Row=7; D='0.5202'; E='  +1.31%  '},
    @{Row=8; D='0.3947'; E='  +3.00%  '},
    @{Row=9; D='0.08360'; E='  +1.51%  '},
    @{Row=10; D='1.135'; E='  +2.23%  '},
    @{Row=11; D='6.303'; E='  +1.78%  '},
    @{Row=12; D='1.918.08'; E='  +3.05%  '},
    @{Row=13; D='20.69'; E='  +0.82%  '},
    @{Row=14; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.038'; E='  +3.38%  '},
    @{Row=15; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.319'; E='  +0.83%  '},
    @{Row=16; D='0.00001113'; E='  +1.43%  '},
    @{Row=17; D='91.86'; E='  +1.52%  '},
    @{Row=18; D='0.06819'; E='  +2.61%  '},
    @{Row=19; D='17.97'; E='  +1.69%  '},
    @{Row=20; D='1.032'; E='  +2.87%  '},
    @{Row=21; D='6.098'; E='  +1.42%  '},
    @{Row=22; D='28.764.86'; E='  +2.64%  '},
    @{Row=23; D='11.27'; E='  +1.87%  '},
    @{Row=24; D='2.274'; E='  +0.80%  '},
    @{Row=25; D='2.133.10'; E='  +2.95%  '},
    @{Row=26; D='162.36'; E='  +3.21%  '},
    @{Row=27; D='21.01'; E='  +2.71%  '},
    @{Row=28; D='2.453'; E='  -2.09%  '},
    @{Row=29; D='127.85'; E='  +2.63%  '},
    @{Row=30; D='0.1061'; E='  -0.39%  '},
    @{Row=31; D='1.056'; E='  +2.43%  '},
    @{Row=32; D='5.989'; E='  +1.49%  '},
    @{Row=33; D='3.686'; E='  +2.60%  '},
    @{Row=34; D='9.468'; E='  +0.47%  '},
    @{Row=35; D='0.02471'; E='  +2.54%  '},
    @{Row=36; D='0.06657'; E='  +2.40%  '},
    @{Row=37; D='0.2226'; E='  +2.30%  '},
    @{Row=38; D='0.6571'; E='  +0.48%  '},
    @{Row=39; D='1.264'; E='  +4.40%  '},
    @{Row=40; D='1.194'; E='  +0.03%  '},
    @{Row=41; D='5.036'; E='  +0.94%  '},
    @{Row=42; D='11.19'; E='  +0.28%  '},
    @{Row=43; D='0.6168'; E='  +0.38%  '},
    @{Row=44; D='13.27'; E='  +2.49%  '},
    @{Row=45; D='3.757'; E='  +2.29%  '},
    @{Row=46; D='1.300'; E='  +1.67%  '},
    @{Row=47; B='EOS'; C='https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; D='1.243'; E='  +2.54%  '},
    @{Row=48; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='2.022'; E='  +0.68%  '},
    @{Row=49; D='122.89'; E='  +1.62%  '},
    @{Row=50; E='  +2.67%  '},
    @{Row=51; D='78.35'; E='  +0.29%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}

# Force the Price column to remain plain text (some values look numeric,
# e.g. '1.036', and Excel would otherwise coerce them to numbers), then
# restore the default (unstyled) cell style so formatting matches the source.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"
foreach ($u in $updates) {
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
}
$priceRange.Style = "Normal"
